# Apply the commit's row updates to the "Artfynd" worksheet.
# The underlying data rows (2, 4, 5 and 6) had their species/observation
# details re-shuffled between records; row 3 and row 7 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 --------------------------------------------------------------
$ws.Range("A2").Value = 111835745
$ws.Range("B2").Value = 77515
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("Q2").Value = 471152.5480076601
$ws.Range("R2").Value = 6810381.652036018
$ws.Range("S2").Value = 5

# --- Row 4 ----------------------------------------------------------------
$ws.Range("A4").Value = 111835718
$ws.Range("Q4").Value = 471101.0270993827
$ws.Range("R4").Value = 6810411.753755242
$ws.Range("S4").Value = 10
$ws.Range("AC4").ClearContents()

# --- Row 5 ------------------------------------------------------------
$ws.Range("A5").Value = 111835838
$ws.Range("B5").Value = 89423
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma"
$ws.Range("H5").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q5").Value = 470914.6782613794
$ws.Range("R5").Value = 6810368.79402096

# --- Row 6 --------------------------------------------------------------
$ws.Range("A6").Value = 111835826
$ws.Range("B6").Value = 56398
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("J6").ClearContents()
$ws.Range("M6").Value = "äldre spår"
$ws.Range("Q6").Value = 470915.776864712
$ws.Range("R6").Value = 6810385.536630718
$ws.Range("AC6").Value = "även hackspettbo, troligen av tret"
$ws.Range("AF6").ClearContents()
